$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 14:13:45"
$ws.Range("A3").Value = "Total filas: 247"
$ws.Range("A24").Value = "03:42:43"
$ws.Range("C24").Value = "14_ABASTO"
$ws.Range("D24").Value = 113
$ws.Range("A25").Value = "04:17:03"
$ws.Range("C25").Value = "215B_EL PATO"
$ws.Range("D25").Value = 78
$ws.Range("C44").Value = "15_ABASTO"
$ws.Range("C45").Value = "23_HERNANDEZ"
$ws.Range("A75").Value = "07:17:57"
$ws.Range("C75").Value = "215B_EL PATO"
$ws.Range("D75").Value = 65
$ws.Range("A76").Value = "07:50:23"
$ws.Range("C76").Value = "16_P MOR-SANTA ANA"
$ws.Range("D76").Value = 32
$ws.Range("A85").Value = "07:17:57"
$ws.Range("C85").Value = "17_ROMERO"
$ws.Range("D85").Value = 96
$ws.Range("A86").Value = "07:50:23"
$ws.Range("C86").Value = "10_OLMOS"
$ws.Range("D86").Value = 63
$ws.Range("C87").Value = "17_ROMERO"
$ws.Range("C88").Value = "225_HARAS DEL SUR"
$ws.Range("C98").Value = "26_HERNANDEZ"
$ws.Range("C99").Value = "23_HERNANDEZ"
$ws.Range("A100").Value = "08:39:38"
$ws.Range("C100").Value = "17_ROMERO"
$ws.Range("D100").Value = 43
$ws.Range("A101").Value = "08:52:26"
$ws.Range("C101").Value = "16_SANTA ANA"
$ws.Range("D101").Value = 30
$ws.Range("A112").Value = "09:28:24"
$ws.Range("C112").Value = "215C_EL PATO"
$ws.Range("D112").Value = 14
$ws.Range("A113").Value = "08:39:38"
$ws.Range("C113").Value = "10_OLMOS"
$ws.Range("D113").Value = 63
$ws.Range("C138").Value = "215C_EL PATO"
$ws.Range("C139").Value = "16_SANTA ANA"
$ws.Range("A169").Value = "11:45:10"
$ws.Range("C169").Value = "23_HERNANDEZ"
$ws.Range("D169").Value = 20
$ws.Range("A170").Value = "12:04:07"
$ws.Range("C170").Value = "16_SANTA ANA"
$ws.Range("D170").Value = 1
$ws.Range("A171").Value = "12:04:07"
$ws.Range("C171").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D171").Value = 2
$ws.Range("A172").Value = "11:00:36"
$ws.Range("C172").Value = "14_ABASTO"
$ws.Range("D172").Value = 66
$ws.Range("C173").Value = "16_P MOR-SANTA ANA"
$ws.Range("C177").Value = "14_ABASTO"
$ws.Range("A178").Value = "10:25:56"
$ws.Range("C178").Value = "26_HERNANDEZ"
$ws.Range("D178").Value = 115
$ws.Range("A179").Value = "11:00:36"
$ws.Range("C179").Value = "215A_EL PATO"
$ws.Range("D179").Value = 80
$ws.Range("C180").Value = "14_ABASTO"
$ws.Range("C182").Value = "215A_EL PATO"
$ws.Range("C196").Value = "15_ABASTO"
$ws.Range("C197").Value = "14_ABASTO"
$ws.Range("A207").Value = "12:04:07"
$ws.Range("C207").Value = "215D_EL PATO"
$ws.Range("D207").Value = 70
$ws.Range("A208").Value = "13:03:48"
$ws.Range("C208").Value = "11_ETCHEVERRY"
$ws.Range("D208").Value = 11
$ws.Range("A210").Value = "13:03:48"
$ws.Range("C210").Value = "26_HERNANDEZ"
$ws.Range("D210").Value = 18
$ws.Range("A211").Value = "12:44:48"
$ws.Range("C211").Value = "10_OLMOS"
$ws.Range("D211").Value = 37
$ws.Range("A224").Value = "13:03:48"
$ws.Range("C224").Value = "16_P MOR-167 Y 521"
$ws.Range("D224").Value = 53
$ws.Range("A225").Value = "12:04:07"
$ws.Range("C225").Value = "225_GOMEZ"
$ws.Range("D225").Value = 112
$ws.Range("A228").Value = "14:13:45"
$ws.Range("B228").Value = "14:13"
$ws.Range("C228").Value = "15_ABASTO"
$ws.Range("D228").Value = 0
$ws.Range("A229").Value = "14:13:45"
$ws.Range("B229").Value = "14:16"
$ws.Range("C229").Value = "27_EL RETIRO"
$ws.Range("D229").Value = 3
$ws.Range("A230").Value = "14:13:45"
$ws.Range("B230").Value = "14:16"
$ws.Range("C230").Value = "14_ABASTO"
$ws.Range("D230").Value = 3
$ws.Range("A231").Value = "14:13:45"
$ws.Range("B231").Value = "14:19"
$ws.Range("C231").Value = "215C_EL PATO"
$ws.Range("D231").Value = 6
$ws.Range("A232").Value = "14:13:45"
$ws.Range("B232").Value = "14:21"
$ws.Range("C232").Value = "26_HERNANDEZ"
$ws.Range("D232").Value = 8
$ws.Range("A233").Value = "14:13:45"
$ws.Range("B233").Value = "14:34"
$ws.Range("C233").Value = "23_HERNANDEZ"
$ws.Range("D233").Value = 21
$ws.Range("A234").Value = "14:13:45"
$ws.Range("B234").Value = "14:44"
$ws.Range("C234").Value = "14_ABASTO"
$ws.Range("D234").Value = 31
$ws.Range("A235").Value = "14:13:45"
$ws.Range("B235").Value = "14:46"
$ws.Range("C235").Value = "16_SANTA ANA"
$ws.Range("D235").Value = 33
$ws.Range("E235").Value = "LP1912"
$ws.Range("A236").Value = "14:13:45"
$ws.Range("B236").Value = "14:51"
$ws.Range("C236").Value = "17_ROMERO"
$ws.Range("D236").Value = 38
$ws.Range("E236").Value = "LP1912"
$ws.Range("A237").Value = "14:13:45"
$ws.Range("B237").Value = "14:56"
$ws.Range("C237").Value = "16_P MOR-SANTA ANA"
$ws.Range("D237").Value = 43
$ws.Range("E237").Value = "LP1912"
$ws.Range("A238").Value = "14:13:45"
$ws.Range("B238").Value = "14:58"
$ws.Range("C238").Value = "215B_EL PATO"
$ws.Range("D238").Value = 45
$ws.Range("E238").Value = "LP1912"
$ws.Range("A239").Value = "14:13:45"
$ws.Range("B239").Value = "15:00"
$ws.Range("C239").Value = "81_EL PELIGRO"
$ws.Range("D239").Value = 47
$ws.Range("E239").Value = "LP1912"
$ws.Range("A240").Value = "14:13:45"
$ws.Range("B240").Value = "15:05"
$ws.Range("C240").Value = "10_OLMOS"
$ws.Range("D240").Value = 52
$ws.Range("E240").Value = "LP1912"
$ws.Range("A241").Value = "14:13:45"
$ws.Range("B241").Value = "15:06"
$ws.Range("C241").Value = "16_SANTA ANA"
$ws.Range("D241").Value = 53
$ws.Range("E241").Value = "LP1912"
$ws.Range("A242").Value = "14:13:45"
$ws.Range("B242").Value = "15:10"
$ws.Range("C242").Value = "17_ROMERO"
$ws.Range("D242").Value = 57
$ws.Range("E242").Value = "LP1912"
$ws.Range("A243").Value = "14:13:45"
$ws.Range("B243").Value = "15:13"
$ws.Range("C243").Value = "11_ETCHEVERRY"
$ws.Range("D243").Value = 60
$ws.Range("E243").Value = "LP1912"
$ws.Range("A244").Value = "14:13:45"
$ws.Range("B244").Value = "15:20"
$ws.Range("C244").Value = "15_ABASTO"
$ws.Range("D244").Value = 67
$ws.Range("E244").Value = "LP1912"
$ws.Range("A245").Value = "14:13:45"
$ws.Range("B245").Value = "15:21"
$ws.Range("C245").Value = "26_HERNANDEZ"
$ws.Range("D245").Value = 68
$ws.Range("E245").Value = "LP1912"
$ws.Range("A246").Value = "14:13:45"
$ws.Range("B246").Value = "15:32"
$ws.Range("C246").Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Range("D246").Value = 79
$ws.Range("E246").Value = "LP1912"
$ws.Range("A247").Value = "14:13:45"
$ws.Range("B247").Value = "15:36"
$ws.Range("C247").Value = "23_HERNANDEZ"
$ws.Range("D247").Value = 83
$ws.Range("E247").Value = "LP1912"
$ws.Range("A248").Value = "14:13:45"
$ws.Range("B248").Value = "15:37"
$ws.Range("C248").Value = "10_OLMOS"
$ws.Range("D248").Value = 84
$ws.Range("E248").Value = "LP1912"
$ws.Range("A249").Value = "14:13:45"
$ws.Range("B249").Value = "15:38"
$ws.Range("C249").Value = "215A_EL PATO"
$ws.Range("D249").Value = 85
$ws.Range("E249").Value = "LP1912"
$ws.Range("A250").Value = "14:13:45"
$ws.Range("B250").Value = "15:46"
$ws.Range("C250").Value = "16_P MOR-167 Y 521"
$ws.Range("D250").Value = 93
$ws.Range("E250").Value = "LP1912"
$ws.Range("A251").Value = "14:13:45"
$ws.Range("B251").Value = "15:53"
$ws.Range("C251").Value = "11_ETCHEVERRY"
$ws.Range("D251").Value = 100
$ws.Range("E251").Value = "LP1912"
$ws.Range("A252").Value = "14:13:45"
$ws.Range("B252").Value = "15:56"
$ws.Range("C252").Value = "27_EL RETIRO"
$ws.Range("D252").Value = 103
$ws.Range("E252").Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 14:13:45"
$ws.Range("A3").Value = "Total filas: 34"
$ws.Range("A37").Value = "14:13:45"
$ws.Range("D37").Value = 6
$ws.Range("A38").Value = "14:13:45"
$ws.Range("D38").Value = 45
$ws.Range("A39").Value = "14:13:45"
$ws.Range("B39").Value = "15:38"
$ws.Range("C39").Value = "215A_EL PATO"
$ws.Range("D39").Value = 85
$ws.Range("E39").Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 14:13:45"
$ws.Range("A3").Value = "Total filas: 38"
$ws.Range("A42").Value = "14:13:45"
$ws.Range("D42").Value = 39
$ws.Range("A43").Value = "14:13:45"
$ws.Range("B43").Value = "15:34"
$ws.Range("C43").Value = "215A_LA PLATA"
$ws.Range("D43").Value = 81
$ws.Range("E43").Value = "L6173"
